$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the dropdown (data validation) that used to live on column B (Ex-Quarry/Delivered)
$ws.Range("B2:B1048576").Validation.Delete()

# Delete the "Transporter Code" / "Transporter Name" columns (C:D); this shifts
# "Customer Code" / "Customer Name" (and everything after them) two columns to the
# left, shrinking the used range from T to R.
$ws.Range("C1:D1").EntireColumn.Delete()

# Add the new "Supplier Code" / "Supplier Name" columns and relabel the old
# "Ex-Quarry/Delivered" column as "Vehicle Weight (KG)".
$ws.Range("E1").Value = "Supplier Code"
$ws.Range("F1").Value = "Supplier Name"
$ws.Range("B1").Value = "Vehicle Weight (KG)"

# The blank E1 cell inherited the "text" number-format style that used to belong to
# the old column G; reset its formatting back to the plain bold header style (copy
# format from A1, which already carries that style).
$ws.Range("A1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Match the saved cursor position recorded in the workbook.
$ws.Range("H14").Select()
